$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be
# auto-converted to a number by Excel (they must stay text, matching the
# original inlineStr cells). Multi-area "A1,B2" Range refs only honour the
# first area here, so set NumberFormat per-cell instead.
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"

$ws.Range('D2').Value = '69.837.68'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '3.529.90'
$ws.Range('E3').Value = '  +0.73%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '605.82'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('D6').Value = '195.50'
$ws.Range('E6').Value = '  +1.95%  '
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D9').Value = '0.203'
$ws.Range('E9').Value = '  -4.75%  '
$ws.Range('D10').Value = '0.646'
$ws.Range('E10').Value = '  -2.47%  '
$ws.Range('D11').Value = '53.49'
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('E12').Value = '  -1.62%  '
$ws.Range('E13').Value = '  -1.74%  '
$ws.Range('D14').Value = '4.090.63'
$ws.Range('E14').Value = '  +0.66%  '
$ws.Range('D15').Value = '594.84'
$ws.Range('E15').Value = '  -3.70%  '
$ws.Range('D16').Value = '69.909.07'
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').Value = '12.71'
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('D18').Value = '18.98'
$ws.Range('E18').Value = '  +0.59%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.541.81'
$ws.Range('E19').Value = '  +1.17%  '
$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D20').Value = '0.123'
$ws.Range('E20').Value = '  +1.88%  '
$ws.Range('D21').Value = '0.984'
$ws.Range('E21').Value = '  -0.70%  '
$ws.Range('D22').Value = '17.79'
$ws.Range('E22').Value = '  -1.16%  '
$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D23').Value = '5.17'
$ws.Range('E23').Value = '  +3.26%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '103.31'
$ws.Range('E24').Value = '  -2.09%  '
$ws.Range('E25').Value = '  -0.52%  '
$ws.Range('D26').Value = '3.05'
$ws.Range('E26').Value = '  +0.60%  '
$ws.Range('D27').Value = '10.81'
$ws.Range('E27').Value = '  -1.77%  '
$ws.Range('D28').Value = '9.52'
$ws.Range('E28').Value = '  -3.87%  '
$ws.Range('D29').Value = '33.23'
$ws.Range('E29').Value = '  -2.79%  '
$ws.Range('D30').Value = '7.05'
$ws.Range('E30').Value = '  -0.95%  '
$ws.Range('D31').Value = '4.23'
$ws.Range('E31').Value = '  +0.95%  '
$ws.Range('D32').Value = '12.33'
$ws.Range('E32').Value = '  -2.64%  '
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').Value = '63.44'
$ws.Range('E34').Value = '  -1.05%  '
$ws.Range('D35').Value = '3.19'
$ws.Range('E35').Value = '  +3.22%  '
$ws.Range('D36').Value = '3.789.74'
$ws.Range('E36').Value = '  +1.35%  '
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '0.0₃0809'
$ws.Range('E37').Value = '  +1.42%  '
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').Value = '512.37'
$ws.Range('E39').Value = '  -2.38%  '
$ws.Range('D40').Value = '0.391'
$ws.Range('E40').Value = '  +0.19%  '
$ws.Range('D41').Value = '3.58'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').Value = '36.46'
$ws.Range('E42').Value = '  -0.96%  '
$ws.Range('D43').Value = '0.134'
$ws.Range('E43').Value = '  -2.90%  '
$ws.Range('E44').Value = '  -3.08%  '
$ws.Range('E45').Value = '  -0.93%  '
$ws.Range('D46').Value = '2.82'
$ws.Range('E46').Value = '  -1.44%  '
$ws.Range('D47').Value = '3.22'
$ws.Range('E47').Value = '  -3.07%  '
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('D49').Value = '8.47'
$ws.Range('E49').Value = '  -3.09%  '
$ws.Range('E50').Value = '  +2.73%  '
$ws.Range('E51').Value = '  +4.98%  '
